$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(according to the population census data)" caption that used
# to sit under the title in A2 - the export no longer includes it.
$ws.Range("A2").ClearContents()

# The blank spacer row (old row 3) is removed entirely, pulling the
# "(sq. km)" row and the data rows up by one.
$ws.Rows("3").Delete()

# Only the 2014 figures are kept now - the 1989 and 2002 columns are
# dropped (old columns B and C), so the former column D (2014 / 997.5)
# becomes the new column B.
$ws.Range("B:C").EntireColumn.Delete()

# The surviving rows are given a uniform, taller row height.
$ws.Rows("1:5").RowHeight = 20.1

# Leave the cursor where it was left in the saved file.
[void]$ws.Range("C9").Select()
